$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two students that were cleaned out of the dataset.
# Deleting the higher row index first keeps the lower index valid.
$ws.Rows.Item(48).Delete()
$ws.Rows.Item(18).Delete()

# Re-apply the refreshed ranking (col A) and roll number (col B) pairing
# for every remaining student row, now re-sorted/re-ranked after cleaning.
$ws.Cells.Item(2, 1).Value = 12
$ws.Cells.Item(2, 2).Value = "G931259509014"
$ws.Cells.Item(3, 1).Value = 20
$ws.Cells.Item(3, 2).Value = "W931252111065"
$ws.Cells.Item(4, 1).Value = 15
$ws.Cells.Item(4, 2).Value = "B931259309004"
$ws.Cells.Item(5, 1).Value = 27
$ws.Cells.Item(5, 2).Value = "G931321113006"
$ws.Cells.Item(6, 1).Value = 31
$ws.Cells.Item(6, 2).Value = "C931321610014"
$ws.Cells.Item(7, 1).Value = 9
$ws.Cells.Item(7, 2).Value = "W931101109061"
$ws.Cells.Item(8, 1).Value = 35
$ws.Cells.Item(8, 2).Value = "D931325309031"
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "D931100608056"
$ws.Cells.Item(10, 1).Value = 10
$ws.Cells.Item(10, 2).Value = "R931325310022"
$ws.Cells.Item(11, 1).Value = 24
$ws.Cells.Item(11, 2).Value = "Q931235212001"
$ws.Cells.Item(12, 1).Value = 18
$ws.Cells.Item(12, 2).Value = "E931100609021"
$ws.Cells.Item(13, 1).Value = 29
$ws.Cells.Item(13, 2).Value = "T887690719015"
$ws.Cells.Item(14, 1).Value = 25
$ws.Cells.Item(14, 2).Value = "R931252710029"
$ws.Cells.Item(15, 1).Value = 50
$ws.Cells.Item(15, 2).Value = "Y931325210014"
$ws.Cells.Item(16, 1).Value = 51
$ws.Cells.Item(16, 2).Value = "Q931101109046"
$ws.Cells.Item(17, 1).Value = 48
$ws.Cells.Item(17, 2).Value = "U931383908031"
$ws.Cells.Item(18, 1).Value = 47
$ws.Cells.Item(18, 2).Value = "M931321110016"
$ws.Cells.Item(19, 1).Value = 52
$ws.Cells.Item(19, 2).Value = "X886463320016"
$ws.Cells.Item(20, 1).Value = 46
$ws.Cells.Item(20, 2).Value = "W931383410018"
$ws.Cells.Item(21, 1).Value = 49
$ws.Cells.Item(21, 2).Value = "X931325208068"
$ws.Cells.Item(22, 1).Value = 56
$ws.Cells.Item(22, 2).Value = "Z931101109005"
$ws.Cells.Item(23, 1).Value = 54
$ws.Cells.Item(23, 2).Value = "K931101109004"
$ws.Cells.Item(24, 1).Value = 44
$ws.Cells.Item(24, 2).Value = "G931383410017"
$ws.Cells.Item(25, 1).Value = 57
$ws.Cells.Item(25, 2).Value = "Y931412017035"
$ws.Cells.Item(26, 1).Value = 58
$ws.Cells.Item(26, 2).Value = "B931235209044"
$ws.Cells.Item(27, 1).Value = 59
$ws.Cells.Item(27, 2).Value = "R931100609009"
$ws.Cells.Item(28, 1).Value = 60
$ws.Cells.Item(28, 2).Value = "N931325209054"
$ws.Cells.Item(29, 1).Value = 61
$ws.Cells.Item(29, 2).Value = "F931252509025"
$ws.Cells.Item(30, 1).Value = 62
$ws.Cells.Item(30, 2).Value = "U931325208066"
$ws.Cells.Item(31, 1).Value = 63
$ws.Cells.Item(31, 2).Value = "M931321009023"
$ws.Cells.Item(32, 1).Value = 53
$ws.Cells.Item(32, 2).Value = "C931100609037"
$ws.Cells.Item(33, 1).Value = 43
$ws.Cells.Item(33, 2).Value = "L931412020028"
$ws.Cells.Item(34, 1).Value = 1
$ws.Cells.Item(34, 2).Value = "Q931258910001"
$ws.Cells.Item(35, 1).Value = 41
$ws.Cells.Item(35, 2).Value = "H931321309010"
$ws.Cells.Item(36, 1).Value = 2
$ws.Cells.Item(36, 2).Value = "L931252709035"
$ws.Cells.Item(37, 1).Value = 4
$ws.Cells.Item(37, 2).Value = "N931101108063"
$ws.Cells.Item(38, 1).Value = 6
$ws.Cells.Item(38, 2).Value = "Q931252109012"
$ws.Cells.Item(39, 1).Value = 8
$ws.Cells.Item(39, 2).Value = "G931383411018"
$ws.Cells.Item(40, 1).Value = 11
$ws.Cells.Item(40, 2).Value = "U931252114001"
$ws.Cells.Item(41, 1).Value = 14
$ws.Cells.Item(41, 2).Value = "M931252110020"
$ws.Cells.Item(42, 1).Value = 16
$ws.Cells.Item(42, 2).Value = "N931253409013"
$ws.Cells.Item(43, 1).Value = 21
$ws.Cells.Item(43, 2).Value = "M931252710007"
$ws.Cells.Item(44, 1).Value = 42
$ws.Cells.Item(44, 2).Value = "M931252916068"
$ws.Cells.Item(45, 1).Value = 22
$ws.Cells.Item(45, 2).Value = "L931412020030"
$ws.Cells.Item(46, 1).Value = 30
$ws.Cells.Item(46, 2).Value = "K931252910051"
$ws.Cells.Item(47, 1).Value = 33
$ws.Cells.Item(47, 2).Value = "C931253116052"
$ws.Cells.Item(48, 1).Value = 34
$ws.Cells.Item(48, 2).Value = "X931412020027"
$ws.Cells.Item(49, 1).Value = 65
$ws.Cells.Item(49, 2).Value = "Q931101008033"
$ws.Cells.Item(50, 1).Value = 36
$ws.Cells.Item(50, 2).Value = "U931101008035"
$ws.Cells.Item(51, 1).Value = 38
$ws.Cells.Item(51, 2).Value = "Q931100609020"
$ws.Cells.Item(52, 1).Value = 39
$ws.Cells.Item(52, 2).Value = "F931100609041"
$ws.Cells.Item(53, 1).Value = 40
$ws.Cells.Item(53, 2).Value = "F931235210018"
$ws.Cells.Item(54, 1).Value = 23
$ws.Cells.Item(54, 2).Value = "D931100609030"
$ws.Cells.Item(55, 1).Value = 66
$ws.Cells.Item(55, 2).Value = "V931101109041"

"edit applied"
